$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 172, shifting all the
# existing data (rows 172:283) down to rows 174:285. This matches the
# dimension change from A1:T283 to A1:T285.
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()

# Row 172 (new): same Mercado/Producto attributes as the rest of the table,
# "Primera" quality, new date + prices + origin.
$ws.Cells.Item(172, 1).Value = 7
$ws.Cells.Item(172, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(172, 3).Value = "Ñuble"
$ws.Cells.Item(172, 4).Value = 45086
$ws.Cells.Item(172, 5).Value = 16
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100101
$ws.Cells.Item(172, 8).Value = "Berries"
$ws.Cells.Item(172, 9).Value = 100101007
$ws.Cells.Item(172, 10).Value = "Kiwi"
$ws.Cells.Item(172, 11).Value = "Hayward"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 60
$ws.Cells.Item(172, 14).Value = 10000
$ws.Cells.Item(172, 15).Value = 10000
$ws.Cells.Item(172, 16).Value = 10000
$ws.Cells.Item(172, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(172, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(172, 19).Value = 556
$ws.Cells.Item(172, 20).Value = 18

# Row 173 (new): same attributes, "Segunda" quality.
$ws.Cells.Item(173, 1).Value = 7
$ws.Cells.Item(173, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(173, 3).Value = "Ñuble"
$ws.Cells.Item(173, 4).Value = 45086
$ws.Cells.Item(173, 5).Value = 16
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100101
$ws.Cells.Item(173, 8).Value = "Berries"
$ws.Cells.Item(173, 9).Value = 100101007
$ws.Cells.Item(173, 10).Value = "Kiwi"
$ws.Cells.Item(173, 11).Value = "Hayward"
$ws.Cells.Item(173, 12).Value = "Segunda"
$ws.Cells.Item(173, 13).Value = 60
$ws.Cells.Item(173, 14).Value = 9000
$ws.Cells.Item(173, 15).Value = 9000
$ws.Cells.Item(173, 16).Value = 9000
$ws.Cells.Item(173, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(173, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(173, 19).Value = 500
$ws.Cells.Item(173, 20).Value = 18

# Make sure the date cells use the same date number format as the rest of
# column D (style index 2 in the original workbook == numFmtId 165).
$ws.Cells.Item(172, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
